$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "NA" values for the duplicate_image_filename column (E)
# for every data row of the first table (rows 2-21).
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
